$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the look of the
# existing header cells (bold, centered, thin border, top-aligned) by
# copying the format from H1 ("IP") before setting the new text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-28.
$data = @(
    @(1,4),
    @(1,4),
    @(7,8),
    @(1,3),
    @(7,8),
    @(9,9),
    @(6,7),
    @(7,8),
    @(1,4),
    @(1,6),
    @(1,5),
    @(1,5),
    @(1,7),
    @(1,7),
    @(1,3),
    @(1,5),
    @(4,7),
    @(4,8),
    @(1,5),
    @(9,9),
    @(7,8),
    @(4,5),
    @(6,6),
    @(8,8),
    @(4,5),
    @(5,6),
    @(3,4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
